$d = $word.ActiveDocument

# 1. Update "Curso (semestre ideal)" line: remove the ", EM (8)" part.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EF (7), EM (8)", $true, $false, $false, $false,
    $false, $true, 1, $false, "Curso (semestre ideal): EF (7)", 2)

# 2. Remove the "Requisitos" heading paragraph and the requirement bullet
#    paragraph that follows it (these were the last two paragraphs of the
#    body, right before the bibliography's trailing section break).
$count = $d.Paragraphs.Count
$req = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Requisitos") {
        $req = $p
        break
    }
}

if ($req -ne $null) {
    $start = $req.Range.Start
    $end = $d.Paragraphs.Item($count).Range.End
    $d.Range($start, $end).Delete()
}
